# Update AD + AE tables in factory_terrain (adds new "factory"/"factory_terrain"
# summary tables in columns H:L, mirroring the existing A:E tables for "factory").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Vertical alignment constant (xlVAlignBottom) used purely to force creation of
# the distinct "data" cell style (numFmtId 0 / no special alignment) that the
# target workbook uses (style index 4) for all the new numeric cells.
$xlVAlignBottom = -4107

function Set-Num($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $value
    $c.VerticalAlignment = $xlVAlignBottom
}

function Set-Str($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

function Set-EmptyStyled($row, $col) {
    $ws.Cells.Item($row, $col).VerticalAlignment = $xlVAlignBottom
}

# ---- row 55: section title ----
Set-Str 55 1 "factory"
Set-Str 55 8 "factory_terrain"

# ---- row 56: "AE" header + factory column labels (H:L) ----
Set-Str 56 8  "AE"
Set-Str 56 9  "factory 1"
Set-Str 56 10 "f2"
Set-Str 56 11 "f3"
Set-Str 56 12 "f4"

# ---- rows 57-59: AE totals for factory_terrain ----
Set-Str 57 8 "total distance"
Set-Num 57 9  24.12
Set-Num 57 10 18.329999999999998
Set-Num 57 11 14.78
Set-Num 57 12 20.61

Set-Str 58 8 "total angle"
Set-Num 58 9  1302.73
Set-Num 58 10 1247.4100000000001
Set-Num 58 11 455.7
Set-Num 58 12 1851.29

Set-Str 59 8 "total time"
Set-Num 59 9  180.68
Set-Num 59 10 135.26
Set-Num 59 11 108.23
Set-Num 59 12 158.47

# ---- row 61: "AE 0.5" header ----
Set-Str 61 8 "AE 0.5"

# ---- rows 62-64: AE 0.5 totals for factory_terrain ----
Set-Str 62 8 "total distance"
Set-Num 62 9  23.56
Set-Num 62 10 18.170000000000002
Set-Num 62 11 14.73
Set-Num 62 12 20.48

Set-Str 63 8 "total angle"
Set-Num 63 9  1146.8399999999999
Set-Num 63 10 916.05
Set-Num 63 11 874.63
Set-Num 63 12 1281.3499999999999

Set-Str 64 8 "total time"
Set-Num 64 9  175.98
Set-Num 64 10 130.74
Set-Num 64 11 115.17
Set-Num 64 12 151.43

# ---- row 67: "AE 10" header ----
Set-Str 67 8 "AE 10"

# ---- rows 68-70: AE 10 totals for factory_terrain ----
Set-Str 68 8 "total distance"
Set-Num 68 9  24.84
Set-Num 68 10 22.97
Set-Num 68 11 14.95
Set-Num 68 12 21.9

Set-Str 69 8 "total angle"
Set-Num 69 9  871.23
Set-Num 69 10 2636.26
Set-Num 69 11 627.74
Set-Num 69 12 2868.05

Set-Str 70 8 "total time"
Set-Num 70 9  176.95
Set-Num 70 10 176.64
Set-Num 70 11 111.96
Set-Num 70 12 179.47

# ---- row 72: spacer row, styled but empty ----
Set-EmptyStyled 72 8
Set-EmptyStyled 72 9
Set-EmptyStyled 72 10
Set-EmptyStyled 72 11

# ---- row 73: section title (second block) ----
Set-Str 73 8 "factory_terrain"

# ---- row 74: "AD" header + factory column labels (H:L) ----
Set-Str 74 8  "AD"
Set-Str 74 9  "factory 1"
Set-Str 74 10 "f2"
Set-Str 74 11 "f3"
Set-Str 74 12 "f4"

# ---- rows 75-77: AD totals for factory_terrain ----
Set-Str 75 8 "total distance"
Set-Num 75 9  23.4
Set-Num 75 10 18.18
Set-Num 75 11 14.21
Set-Num 75 12 20.49

Set-Str 76 8 "total angle"
Set-Num 76 9  811.91
Set-Num 76 10 932.45
Set-Num 76 11 2126.2600000000002
Set-Num 76 12 1131.79

Set-Str 77 8 "total time"
Set-Num 77 9  170.24
Set-Num 77 10 130.53
Set-Num 77 11 155.55000000000001
Set-Num 77 12 152.51

# ---- row 79: "AD 0.5" header ----
Set-Str 79 8 "AD 0.5"

# ---- rows 80-82: AD 0.5 totals for factory_terrain ----
Set-Str 80 8 "total distance"
Set-Num 80 9  23.36
Set-Num 80 10 18.18
Set-Num 80 11 14.69
Set-Num 80 12 20.48

Set-Str 81 8 "total angle"
Set-Num 81 9  1081.4000000000001
Set-Num 81 10 920.26
Set-Num 81 11 1140.42
Set-Num 81 12 1300.8399999999999

Set-Str 82 8 "total time"
Set-Num 82 9  173.5
Set-Num 82 10 129.93
Set-Num 82 11 113.74
Set-Num 82 12 153.76

# ---- row 85: "AD 10" header ----
Set-Str 85 8 "AD 10"

# ---- rows 86-88: AD 10 totals for factory_terrain ----
Set-Str 86 8 "total distance"
Set-Num 86 9  24.01
Set-Num 86 10 18.809999999999999
Set-Num 86 11 14.95
Set-Num 86 12 24.82

Set-Str 87 8 "total angle"
Set-Num 87 9  542.54
Set-Num 87 10 951.41
Set-Num 87 11 627.74
Set-Num 87 12 1249.6600000000001

Set-Str 88 8 "total time"
Set-Num 88 9  176.78
Set-Num 88 10 136.79
Set-Num 88 11 111.96
Set-Num 88 12 187.91

# ---- update the active selection to mirror the author's final cursor position ----
$ws.Range("I75").Select()

Write-Output "edit complete"
